$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 94.7
$ws.Cells.Item(9, 9).Value = 99.55556
$ws.Cells.Item(9, 11).Value = 99.55556
$ws.Cells.Item(9, 13).Value = 69.44444

$ws.Cells.Item(18, 8).Value = 286.14285
$ws.Cells.Item(18, 9).Value = 286.14285
$ws.Cells.Item(18, 11).Value = 286.14285
$ws.Cells.Item(18, 13).Value = -2.14285000000001

$ws.Cells.Item(33, 8).Value = 3380422
$ws.Cells.Item(33, 9).Value = 4504896
$ws.Cells.Item(33, 11).Value = 4504896
$ws.Cells.Item(33, 13).Value = -4504667

$ws.Cells.Item(118, 8).Value = 573.125
$ws.Cells.Item(118, 9).Value = 573.125
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 1719.375
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = -62.375
$ws.Cells.Item(118, 14).ClearContents()

$ws.Cells.Item(129, 8).Value = 8131.4443
$ws.Cells.Item(129, 10).Value = 9099.25
$ws.Cells.Item(129, 12).Value = 27297.75
$ws.Cells.Item(129, 14).Value = -37297.75

$ws.Cells.Item(132, 8).Value = 1729.4036
$ws.Cells.Item(132, 9).Value = 1755.4231
$ws.Cells.Item(132, 11).Value = 5266.2693
$ws.Cells.Item(132, 13).Value = -2736.2693

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1123.3334
$ws.Cells.Item(110, 9).Value = 1153.8667
$ws.Cells.Item(110, 10).Value = 1085.1666
$ws.Cells.Item(110, 11).Value = 1153.8667
$ws.Cells.Item(110, 12).Value = 1085.1666
$ws.Cells.Item(110, 13).Value = 891.1333
$ws.Cells.Item(110, 14).Value = -5175.1666

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 11642.391
$ws.Cells.Item(20, 9).Value = 13084.311
$ws.Cells.Item(20, 10).Value = 8157.75
$ws.Cells.Item(20, 11).Value = 13084.311
$ws.Cells.Item(20, 12).Value = 8157.75
$ws.Cells.Item(20, 13).Value = -12837.311
$ws.Cells.Item(20, 14).Value = -8651.75

$ws.Cells.Item(94, 8).Value = 1408.25
$ws.Cells.Item(94, 9).Value = 1271.2354
$ws.Cells.Item(94, 11).Value = 1271.2354
$ws.Cells.Item(94, 13).Value = -820.2354

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2213.12
$ws.Cells.Item(16, 9).Value = 2292.9167
$ws.Cells.Item(16, 10).Value = 298
$ws.Cells.Item(16, 11).Value = 2292.9167
$ws.Cells.Item(16, 12).Value = 298
$ws.Cells.Item(16, 13).Value = -2005.9167
$ws.Cells.Item(16, 14).Value = -872

$ws.Cells.Item(58, 8).Value = 3455.0908
$ws.Cells.Item(58, 9).Value = 2126.7273
$ws.Cells.Item(58, 10).Value = 4783.4546
$ws.Cells.Item(58, 11).Value = 2126.7273
$ws.Cells.Item(58, 12).Value = 4783.4546
$ws.Cells.Item(58, 13).Value = -1923.7273
$ws.Cells.Item(58, 14).Value = -5189.4546

$ws.Cells.Item(99, 8).Value = 10039.548
$ws.Cells.Item(99, 9).Value = 5765.846
$ws.Cells.Item(99, 10).Value = 11955.345
$ws.Cells.Item(99, 11).Value = 5765.846
$ws.Cells.Item(99, 12).Value = 11955.345
$ws.Cells.Item(99, 13).Value = -4267.846
$ws.Cells.Item(99, 14).Value = -14951.345

$ws.Cells.Item(113, 8).Value = 2213.12
$ws.Cells.Item(113, 9).Value = 2292.9167
$ws.Cells.Item(113, 10).Value = 298
$ws.Cells.Item(113, 11).Value = 2292.9167
$ws.Cells.Item(113, 12).Value = 298
$ws.Cells.Item(113, 13).Value = -122.9167000000002
$ws.Cells.Item(113, 14).Value = -4638

$ws.Cells.Item(126, 8).Value = 10039.548
$ws.Cells.Item(126, 9).Value = 5765.846
$ws.Cells.Item(126, 10).Value = 11955.345
$ws.Cells.Item(126, 11).Value = 17297.538
$ws.Cells.Item(126, 12).Value = 35866.035
$ws.Cells.Item(126, 13).Value = -14827.538
$ws.Cells.Item(126, 14).Value = -40806.035

$ws.Cells.Item(136, 8).Value = 3455.0908
$ws.Cells.Item(136, 9).Value = 2126.7273
$ws.Cells.Item(136, 10).Value = 4783.4546
$ws.Cells.Item(136, 11).Value = 6380.1819
$ws.Cells.Item(136, 12).Value = 14350.3638
$ws.Cells.Item(136, 13).Value = -3830.1819
$ws.Cells.Item(136, 14).Value = -19450.3638

$ws.Cells.Item(141, 8).Value = 259444.9
$ws.Cells.Item(141, 10).Value = 306181.12
$ws.Cells.Item(141, 12).Value = 306181.12
$ws.Cells.Item(141, 14).Value = -316541.12

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 3335
$ws.Cells.Item(14, 9).Value = 3335
$ws.Cells.Item(14, 11).Value = 10005
$ws.Cells.Item(14, 13).Value = -9832

$ws.Cells.Item(121, 8).Value = 2454631.2
$ws.Cells.Item(121, 9).Value = 2709.889
$ws.Cells.Item(121, 10).Value = 5849599
$ws.Cells.Item(121, 11).Value = 8129.667
$ws.Cells.Item(121, 12).Value = 17548797
$ws.Cells.Item(121, 13).Value = -6819.667
$ws.Cells.Item(121, 14).Value = -17551417

$ws.Cells.Item(123, 8).Value = 1856.75
$ws.Cells.Item(123, 9).Value = 476
$ws.Cells.Item(123, 10).Value = 5999
$ws.Cells.Item(123, 11).Value = 1428
$ws.Cells.Item(123, 12).Value = 17997
$ws.Cells.Item(123, 13).Value = 1022
$ws.Cells.Item(123, 14).Value = -22897

$ws.Cells.Item(136, 8).Value = 4000
$ws.Cells.Item(136, 9).Value = 4000
$ws.Cells.Item(136, 11).Value = 12000
$ws.Cells.Item(136, 13).Value = -6900

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 25590
$ws.Cells.Item(34, 10).Value = 25590
$ws.Cells.Item(34, 12).Value = 25590
$ws.Cells.Item(34, 14).Value = -26126

$ws.Cells.Item(76, 8).Value = 25590
$ws.Cells.Item(76, 10).Value = 25590
$ws.Cells.Item(76, 12).Value = 25590
$ws.Cells.Item(76, 14).Value = -26220

$ws.Cells.Item(79, 8).Value = 25590
$ws.Cells.Item(79, 10).Value = 25590
$ws.Cells.Item(79, 12).Value = 25590
$ws.Cells.Item(79, 14).Value = -27774

$ws.Cells.Item(102, 8).Value = 1956.4615
$ws.Cells.Item(102, 9).Value = 1874.9166
$ws.Cells.Item(102, 11).Value = 1874.9166
$ws.Cells.Item(102, 13).Value = -252.9166

$ws.Cells.Item(122, 8).Value = 4166.3335
$ws.Cells.Item(122, 9).Value = 4166.3335
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 12499.0005
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -10049.0005
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 15969.05
$ws.Cells.Item(132, 9).Value = 19782.166
$ws.Cells.Item(132, 10).Value = 10249.375
$ws.Cells.Item(132, 11).Value = 59346.49800000001
$ws.Cells.Item(132, 12).Value = 30748.125
$ws.Cells.Item(132, 13).Value = -56816.49800000001
$ws.Cells.Item(132, 14).Value = -35808.125

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1382.2222
$ws.Cells.Item(22, 9).Value = 887.8
$ws.Cells.Item(22, 11).Value = 887.8
$ws.Cells.Item(22, 13).Value = -592.8

$ws.Cells.Item(27, 8).Value = 1382.2222
$ws.Cells.Item(27, 9).Value = 887.8
$ws.Cells.Item(27, 11).Value = 887.8
$ws.Cells.Item(27, 13).Value = -780.8

$ws.Cells.Item(40, 8).Value = 5767.375
$ws.Cells.Item(40, 9).Value = 5767.375
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 5767.375
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -5631.375
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(43, 8).Value = 89856.71000000001
$ws.Cells.Item(43, 10).Value = 92799.60000000001
$ws.Cells.Item(43, 12).Value = 92799.60000000001
$ws.Cells.Item(43, 14).Value = -93185.60000000001

$ws.Cells.Item(61, 8).Value = 790
$ws.Cells.Item(61, 9).Value = 737.38464
$ws.Cells.Item(61, 11).Value = 737.38464
$ws.Cells.Item(61, 13).Value = -535.38464

$ws.Cells.Item(113, 8).Value = 790
$ws.Cells.Item(113, 9).Value = 737.38464
$ws.Cells.Item(113, 11).Value = 737.38464
$ws.Cells.Item(113, 13).Value = 1432.61536

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()
